$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: Machine Group list (already alphabetically ordered) ---
$groups = @(
    "ASSY",
    "BOMAR",
    "BOMAR4",
    "BOOMAR",
    "E50",
    "ITC",
    "PCM_GROUP_1",
    "PLASMA",
    "PTC",
    "SCM_GROUP_1",
    "SHEARING",
    "SLOT COMMONISATION",
    "SMC_GROUP_1",
    "SMS_GROUP_1",
    "SOCO",
    "TBV_GROUP_1",
    "TCA_GROUP_1",
    "TCM_GROUP_1",
    "TRUMPF"
)

$row = 2
foreach ($g in $groups) {
    $ws.Cells.Item($row, 1).Value = $g
    $row++
}

# --- Column B: Machine Type header + values, entered grouped by type ---
$ws.Range("B1").Value = "Machine Type"

$ws.Cells.Item(3, 2).Value = "Band saw"
$ws.Cells.Item(4, 2).Value = "Band saw"
$ws.Cells.Item(5, 2).Value = "Band saw"

$ws.Cells.Item(8, 2).Value = "CNC"
$ws.Cells.Item(9, 2).Value = "CNC"
$ws.Cells.Item(10, 2).Value = "CNC"
$ws.Cells.Item(11, 2).Value = "CNC"
$ws.Cells.Item(12, 2).Value = "CNC"

$ws.Cells.Item(2, 2).Value = "Manual"
$ws.Cells.Item(6, 2).Value = "Manual"
$ws.Cells.Item(7, 2).Value = "Manual"
$ws.Cells.Item(13, 2).Value = "Manual"
$ws.Cells.Item(14, 2).Value = "Manual"
$ws.Cells.Item(15, 2).Value = "Manual"
$ws.Cells.Item(16, 2).Value = "Manual"
$ws.Cells.Item(17, 2).Value = "Manual"
$ws.Cells.Item(18, 2).Value = "Manual"
$ws.Cells.Item(19, 2).Value = "Manual"
$ws.Cells.Item(20, 2).Value = "Manual"

# --- Column A width (best-fit to content) ---
$ws.Columns.Item(1).AutoFit()

# --- View / selection state ---
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("C3").Select()
